$wb = $excel.ActiveWorkbook

# A new handoff report was generated for the b98b5258 file entry,
# updating its "Latest Handoff Date/Datetime" on the Overview sheet
# and on each per-language handoff detail sheet.

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D4").Value = "2016-19-11 08:19:49"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4").Value = "2016-03-11 08:19:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4").Value = "2016-03-11 08:19:49"
